$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Review Status" value: Open -> Close
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Open", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Close", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: write text into a (currently empty) table cell as a single run.
# ---------------------------------------------------------------------------
function Set-CellText($cell, [string]$text) {
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# Helper: write text into an (empty) cell as two separate runs split at
# $splitLen characters, with no leftover visible formatting difference
# between the two runs (relies on a transient Bold flip to force the
# engine to materialize a run boundary, then flips it back off).
function Set-CellTextTwoRuns($cell, [string]$part1, [string]$part2) {
    $r = $cell.Range
    $r.End = $r.End - 1
    $start = $r.Start
    $r.Text = $part1 + $part2
    $boundary = $start + $part1.Length
    $tail = $d.Range($boundary, $boundary + $part2.Length)
    $tail.Font.Bold = 1
    $tail.Font.Bold = 0
}

# Helper: write text into an (empty) cell as two runs separated by a
# (bookmarked) "_GoBack" bookmark, matching Word's own last-edit marker.
function Set-CellTextWithGoBack($cell, [string]$part1, [string]$part2) {
    $r = $cell.Range
    $r.End = $r.End - 1
    $start = $r.Start
    $r.Text = $part1 + $part2
    $boundary = $start + $part1.Length
    $insPoint = $d.Range($boundary, $boundary)
    $d.Bookmarks.Add("_GoBack", $insPoint) | Out-Null
}

$t = $d.Tables.Item(2)

# ---------------------------------------------------------------------------
# 2) Row "1" (Document Section 5.6 / Blue comment present): fill in the
#    "Responsible person/Planned date for completion" and
#    "Completion(Name/Date)" cells.
# ---------------------------------------------------------------------------
Set-CellTextTwoRuns $t.Cell(2, 5) "Rubén Cocoletzi" " 22-Jan-21"
Set-CellTextWithGoBack $t.Cell(2, 6) "Removing blue c" "omments"

# ---------------------------------------------------------------------------
# 3) Row "2" (Document section 8 / Blue comment present): fill in the
#    classification, responsible person and completion cells.
# ---------------------------------------------------------------------------
Set-CellText $t.Cell(3, 4) "R"
Set-CellTextTwoRuns $t.Cell(3, 5) "Rubén Cocoletzi" " 22-Jan-21"
Set-CellText $t.Cell(3, 6) "Removing blue comments"
